# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets.
# The data is an automated refresh of crawled counts (gh-pages output),
# so only the F-column numeric values change; everything else is untouched.

$wb = $excel.ActiveWorkbook

# Row -> new F-column value, for worksheet "展览"
$exhibitionUpdates = @{
    2  = 2527
    3  = 340
    5  = 1431
    6  = 1113
    7  = 320
    11 = 106
    13 = 8787
    14 = 379
    16 = 267
    20 = 605
    24 = 2046
    25 = 2113
    27 = 1801
    32 = 58
    33 = 113
    34 = 194
    35 = 13
    38 = 262
    40 = 759
    42 = 271
}

# Row -> new F-column value, for worksheet "全部类型"
$allTypesUpdates = @{
    2  = 2527
    3  = 340
    5  = 1431
    7  = 1113
    8  = 320
    12 = 106
    14 = 8787
    15 = 379
    18 = 267
    22 = 605
    26 = 2046
    27 = 2113
    29 = 1801
    34 = 58
    35 = 113
    36 = 194
    37 = 13
    40 = 262
    46 = 759
    49 = 271
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
